$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 4 (item id 5470)
$ws.Range("H4").Value = 983.6667
$ws.Range("I4").Value = 983.6667
$ws.Range("K4").Value = 983.6667
$ws.Range("M4").Value = -869.6667
# Row 9 (item id 5487)
$ws.Range("H9").Value = 169.15384
$ws.Range("I9").Value = 92.625
$ws.Range("K9").Value = 92.625
$ws.Range("M9").Value = 76.375
# Row 18 (item id 5471)
$ws.Range("H18").Value = 1986.1818
$ws.Range("I18").Value = 1986.1818
$ws.Range("K18").Value = 1986.1818
$ws.Range("M18").Value = -1702.1818
# Row 19 (item id 7015)
$ws.Range("H19").Value = 933
$ws.Range("J19").Value = 789.5
$ws.Range("L19").Value = 789.5
$ws.Range("N19").Value = -1139.5
# Row 28 (item id 27772)
$ws.Range("H28").Value = 942.1539
$ws.Range("I28").Value = 767.4545000000001
$ws.Range("J28").Value = 1903
$ws.Range("K28").Value = 767.4545000000001
$ws.Range("L28").Value = 1903
$ws.Range("M28").Value = -282.4545000000001
$ws.Range("N28").Value = -2873
# Row 92 (item id 19901)
$ws.Range("H92").Value = 1139.1
$ws.Range("I92").Value = 1154.5555
$ws.Range("K92").Value = 1154.5555
$ws.Range("M92").Value = 93.44450000000006
# Row 96 (item id 19894)
$ws.Range("H96").Value = 2895.6
$ws.Range("I96").Value = 2122.2856
$ws.Range("J96").Value = 4700
$ws.Range("K96").Value = 6366.8568
$ws.Range("L96").Value = 14100
$ws.Range("M96").Value = -4993.8568
$ws.Range("N96").Value = -16846
# Row 101 (item id 19884)
$ws.Range("H101").Value = 400
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 400
$ws.Range("K101").Value = 0
$ws.Range("N101").Value = -4444
$ws.Range("M101").ClearContents()
# Row 116 (item id 27778)
$ws.Range("H116").Value = 3971.2856
$ws.Range("I116").Value = 3533
$ws.Range("J116").Value = 4300
$ws.Range("K116").Value = 3533
$ws.Range("L116").Value = 4300
$ws.Range("M116").Value = -91
$ws.Range("N116").Value = -11184
# Row 135 (item id 44047)
$ws.Range("H135").Value = 1130.8823
$ws.Range("I135").Value = 1021.11536
$ws.Range("K135").Value = 9190.03824
$ws.Range("M135").Value = -6655.03824
# Row 137 (item id 44013)
$ws.Range("H137").Value = 1756.7693
$ws.Range("I137").Value = 1384.0646
$ws.Range("K137").Value = 4152.1938
$ws.Range("M137").Value = -1602.1938
# Row 138 (item id 44169)
$ws.Range("H138").Value = 3609.9048
$ws.Range("J138").Value = 2618.889
$ws.Range("L138").Value = 7856.667
$ws.Range("N138").Value = -18136.667

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (item id 44147)
$ws.Range("H32").Value = 5524.485
$ws.Range("I32").Value = 3812.1853
$ws.Range("J32").Value = 13229.833
$ws.Range("K32").Value = 3812.1853
$ws.Range("L32").Value = 13229.833
$ws.Range("M32").Value = -3525.1853
$ws.Range("N32").Value = -13803.833
# Row 43 (item id 21715)
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("N43").ClearContents()
# Row 61 (item id 43999)
$ws.Range("H61").Value = 1337.138
$ws.Range("I61").Value = 1236.2593
$ws.Range("K61").Value = 1236.2593
$ws.Range("M61").Value = -1024.2593
# Row 74 (item id 44000)
$ws.Range("H74").Value = 1681.4166
$ws.Range("I74").Value = 1048.381
$ws.Range("K74").Value = 1048.381
$ws.Range("M74").Value = -174.3810000000001
# Row 77 (item id 44000)
$ws.Range("H77").Value = 1681.4166
$ws.Range("I77").Value = 1048.381
$ws.Range("K77").Value = 5241.905000000001
$ws.Range("M77").Value = -873.9050000000007
# Row 102 (item id 19945)
$ws.Range("H102").Value = 1468
$ws.Range("J102").Value = 1814.6666
$ws.Range("L102").Value = 1814.6666
$ws.Range("N102").Value = -5058.6666
# Row 136 (item id 43999)
$ws.Range("H136").Value = 1337.138
$ws.Range("I136").Value = 1236.2593
$ws.Range("K136").Value = 3708.7779
$ws.Range("M136").Value = -1158.7779

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20 (item id 14149)
$ws.Range("H20").Value = 3572.8635
$ws.Range("I20").Value = 3802.4285
$ws.Range("K20").Value = 3802.4285
$ws.Range("M20").Value = -3555.4285
# Row 134 (item id 43998)
$ws.Range("H134").Value = 2813
$ws.Range("I134").Value = 2401.5833
$ws.Range("K134").Value = 7204.749899999999
$ws.Range("M134").Value = -4669.749899999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22 (item id 5367)
$ws.Range("H22").Value = 596.6
$ws.Range("J22").Value = 595
$ws.Range("L22").Value = 595
$ws.Range("N22").Value = -1295
# Row 31 (item id 44023)
$ws.Range("H31").Value = 6411.0625
$ws.Range("I31").Value = 5837.5713
$ws.Range("J31").Value = 6857.1113
$ws.Range("K31").Value = 5837.5713
$ws.Range("L31").Value = 6857.1113
$ws.Range("M31").Value = -5542.5713
$ws.Range("N31").Value = -7447.1113
# Row 34 (item id 44023)
$ws.Range("H34").Value = 6411.0625
$ws.Range("I34").Value = 5837.5713
$ws.Range("J34").Value = 6857.1113
$ws.Range("K34").Value = 5837.5713
$ws.Range("L34").Value = 6857.1113
$ws.Range("M34").Value = -5635.5713
$ws.Range("N34").Value = -7261.1113
# Row 132 (item id 44019)
$ws.Range("H132").Value = 2864.5908
$ws.Range("I132").Value = 2524.8572
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 7574.571599999999
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -5044.571599999999
$ws.Range("N132").Value = -35057
# Row 134 (item id 44020)
$ws.Range("H134").Value = 3335.3635
$ws.Range("I134").Value = 2185
$ws.Range("J134").Value = 4715.8
$ws.Range("K134").Value = 6555
$ws.Range("L134").Value = 14147.4
$ws.Range("M134").Value = -4020
$ws.Range("N134").Value = -19217.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4 (item id 4650)
$ws.Range("H4").Value = 28279496
$ws.Range("I4").Value = 33420880
$ws.Range("J4").Value = 1893.1666
$ws.Range("K4").Value = 100262640
$ws.Range("L4").Value = 5679.4998
$ws.Range("M4").Value = -100262528
$ws.Range("N4").Value = -5903.4998
# Row 5 (item id 43974)
$ws.Range("H5").Value = 423.36365
$ws.Range("I5").Value = 382.1111
$ws.Range("J5").Value = 609
$ws.Range("K5").Value = 1146.3333
$ws.Range("L5").Value = 1827
$ws.Range("M5").Value = -1034.3333
$ws.Range("N5").Value = -2051
# Row 87 (item id 12864)
$ws.Range("H87").Value = 800
$ws.Range("J87").Value = 2000
$ws.Range("L87").Value = 6000
$ws.Range("N87").Value = -8496
# Row 90 (item id 12864)
$ws.Range("H90").Value = 800
$ws.Range("J90").Value = 2000
$ws.Range("L90").Value = 18000
$ws.Range("N90").Value = -30480
# Row 117 (item id 27870)
$ws.Range("H117").Value = 3437.6667
$ws.Range("J117").Value = 3992
$ws.Range("L117").Value = 11976
$ws.Range("N117").Value = -18860
# Row 121 (item id 27878)
$ws.Range("H121").Value = 1317.5714
$ws.Range("J121").Value = 1660.6
$ws.Range("L121").Value = 4981.799999999999
$ws.Range("N121").Value = -7601.799999999999
# Row 135 (item id 43974)
$ws.Range("H135").Value = 423.36365
$ws.Range("I135").Value = 382.1111
$ws.Range("J135").Value = 609
$ws.Range("K135").Value = 3438.9999
$ws.Range("L135").Value = 5481
$ws.Range("M135").Value = -903.9999000000003
$ws.Range("N135").Value = -10551

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61 (item id 27740)
$ws.Range("H61").Value = 4101.231
$ws.Range("I61").Value = 4239.4165
$ws.Range("J61").Value = 2443
$ws.Range("K61").Value = 4239.4165
$ws.Range("L61").Value = 2443
$ws.Range("M61").Value = -4037.4165
$ws.Range("N61").Value = -2847
# Row 113 (item id 27740)
$ws.Range("H113").Value = 4101.231
$ws.Range("I113").Value = 4239.4165
$ws.Range("J113").Value = 2443
$ws.Range("K113").Value = 4239.4165
$ws.Range("L113").Value = 2443
$ws.Range("M113").Value = -2069.4165
$ws.Range("N113").Value = -6783
# Row 132 (item id 44058)
$ws.Range("H132").Value = 4264
$ws.Range("I132").Value = 3687.8096
$ws.Range("J132").Value = 5128.2856
$ws.Range("K132").Value = 11063.4288
$ws.Range("L132").Value = 15384.8568
$ws.Range("M132").Value = -8533.4288
$ws.Range("N132").Value = -20444.8568
# Row 136 (item id 44060)
$ws.Range("H136").Value = 2760.35
$ws.Range("I136").Value = 2663.5789
$ws.Range("K136").Value = 7990.736699999999
$ws.Range("M136").Value = -5440.736699999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136 (item id 44031)
$ws.Range("H136").Value = 2704.48
$ws.Range("I136").Value = 1268.2051
$ws.Range("J136").Value = 7796.727
$ws.Range("K136").Value = 3804.615299999999
$ws.Range("L136").Value = 23390.181
$ws.Range("M136").Value = -1254.615299999999
$ws.Range("N136").Value = -28490.181
